$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.716.04"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "3.383.33"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "575.98"
$ws.Range("E5").Value = "  -3.83%  "
$ws.Range("D6").Value = "133.78"
$ws.Range("E6").Value = "  -6.30%  "
$ws.Range("D8").Value = "3.377.12"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("E10").Value = "  -10.39%  "
$ws.Range("E11").Value = "  -10.31%  "
$ws.Range("E12").Value = "  -8.19%  "
$ws.Range("D13").Value = "3.958.33"
$ws.Range("E13").Value = "  -3.67%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "3.418.85"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  -11.40%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "0.115"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "64.704.82"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "25.82"
$ws.Range("E18").Value = "  -8.79%  "
$ws.Range("E19").Value = "  -14.24%  "
$ws.Range("D20").Value = "5.74"
$ws.Range("E20").Value = "  -7.09%  "
$ws.Range("D21").Value = "13.36"
$ws.Range("E21").Value = "  -6.00%  "
$ws.Range("D22").Value = "377.10"
$ws.Range("E22").Value = "  -8.93%  "
$ws.Range("D23").Value = "0.544"
$ws.Range("E23").Value = "  -8.54%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "71.51"
$ws.Range("E25").Value = "  -7.49%  "
$ws.Range("D26").Value = "3.521.00"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("E27").Value = "  -11.04%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "6.91"
$ws.Range("E29").Value = "  -9.82%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.16"
$ws.Range("E30").Value = "  -11.04%  "
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  -10.40%  "
$ws.Range("D32").Value = "3.395.21"
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -7.17%  "
$ws.Range("D35").Value = "22.64"
$ws.Range("E35").Value = "  -6.58%  "
$ws.Range("D36").Value = "169.33"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("D37").Value = "6.55"
$ws.Range("E37").Value = "  -12.33%  "
$ws.Range("E38").Value = "  -12.60%  "
$ws.Range("D39").Value = "1.43"
$ws.Range("E39").Value = "  -8.02%  "
$ws.Range("D40").Value = "4.59"
$ws.Range("E40").Value = "  -12.35%  "
$ws.Range("D41").Value = "0.0741"
$ws.Range("E41").Value = "  -8.90%  "
$ws.Range("D42").Value = "0.803"
$ws.Range("E42").Value = "  -5.93%  "
$ws.Range("D43").Value = "42.81"
$ws.Range("E43").Value = "  -5.39%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "4.27"
$ws.Range("E45").Value = "  -15.20%  "
$ws.Range("E46").Value = "  -10.97%  "
$ws.Range("D47").Value = "1.07"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "21.56"
$ws.Range("E48").Value = "  -5.95%  "
$ws.Range("D49").Value = "6.37"
$ws.Range("E49").Value = "  -8.74%  "
$ws.Range("D50").Value = "2.140.32"
$ws.Range("E50").Value = "  -8.89%  "
$ws.Range("E51").Value = "  -15.41%  "
